$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.03326466666666666
$ws.Range("H2").Value = 0.09979399999999999
$ws.Range("M2").Value = 101.8783343333333
$ws.Range("N2").Value = 305.635003
$ws.Range("O2").Value = 0.9137375742483709
$ws.Range("P2").Value = 0.913737574248371
$ws.Range("Q2").Value = 3.388948832153555
$ws.Range("R2").Value = 30.500539489382
$ws.Range("S2").Value = 0.9137375742483709
$ws.Range("T2").Value = 0.913737574248371

# Row 3
$ws.Range("G3").Value = 0.03326466666666666
$ws.Range("H3").Value = 0.09979399999999999
$ws.Range("M3").Value = 0.050239
$ws.Range("O3").Value = 0.0004505890510780002
$ws.Range("P3").Value = 0.0004505890510780002
$ws.Range("S3").Value = 0.0004505890510780002
$ws.Range("T3").Value = 0.0004505890510780002

# Row 4
$ws.Range("G4").Value = 0.03326466666666666
$ws.Range("H4").Value = 0.09979399999999999
$ws.Range("M4").Value = 3.444232
$ws.Range("N4").Value = 10.332696
$ws.Range("O4").Value = 0.03089100556484967
$ws.Range("P4").Value = 0.03089100556484967
$ws.Range("Q4").Value = 0.1145712294026667
$ws.Range("R4").Value = 1.031141064624
$ws.Range("S4").Value = 0.03089100556484967
$ws.Range("T4").Value = 0.03089100556484967

# Row 5
$ws.Range("G5").Value = 0.03326466666666666
$ws.Range("H5").Value = 0.09979399999999999
$ws.Range("M5").Value = 0.07069366666666667
$ws.Range("N5").Value = 0.212081
$ws.Range("O5").Value = 0.0006340451079949399
$ws.Range("P5").Value = 0.00063404510799494
$ws.Range("Q5").Value = 0.002351601257111111
$ws.Range("R5").Value = 0.021164411314
$ws.Range("S5").Value = 0.0006340451079949399
$ws.Range("T5").Value = 0.00063404510799494

# Row 6
$ws.Range("G6").Value = 0.03326466666666666
$ws.Range("H6").Value = 0.09979399999999999
$ws.Range("M6").Value = 1.402613333333333
$ws.Range("N6").Value = 4.20784
$ws.Range("O6").Value = 0.01257991223742545
$ws.Range("P6").Value = 0.01257991223742545
$ws.Range("Q6").Value = 0.04665746499555555
$ws.Range("R6").Value = 0.41991718496
$ws.Range("S6").Value = 0.01257991223742545
$ws.Range("T6").Value = 0.01257991223742545

# Row 7
$ws.Range("G7").Value = 0.03326466666666666
$ws.Range("H7").Value = 0.09979399999999999
$ws.Range("M7").Value = 4.650161000000001
$ws.Range("N7").Value = 13.950483
$ws.Range("O7").Value = 0.04170687379028095
$ws.Range("P7").Value = 0.04170687379028095
$ws.Range("Q7").Value = 0.1546860556113333
$ws.Range("R7").Value = 1.392174500502
$ws.Range("S7").Value = 0.04170687379028095
$ws.Range("T7").Value = 0.04170687379028095
